# Atualização de bases das ligas, do dia: 23-02-2024 às 08:18
# Applies corrected/updated match data for "Australia ALeague" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        [hashtable]$Values
    )
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

# Row 68 (id 66) — corrected match data (was swapped with row 69)
Set-Row 68 @{
    'B' = 5404714
    'F' = 'Newcastle Jets'
    'G' = 'Perth Glory'
    'H' = 2
    'I' = 2
    'J' = 'D'
    'L' = 3.3
    'M' = 3
    'N' = 2.15
    'O' = 3.75
    'P' = 3.2
    'R' = 1.9
    'S' = 1.95
    'U' = 1.9
    'V' = 1.95
    'W' = -1
    'X' = 2.75
    'Z' = -0.5
    'AA' = 0.475
    'AB' = 0.8999999999999999
    'AC' = -1
}

# Row 69 (id 67) — corrected match data (was swapped with row 68)
Set-Row 69 @{
    'B' = 5404713
    'F' = 'Brisbane Roar'
    'G' = 'Western United FC'
    'H' = 1
    'I' = 0
    'J' = 'H'
    'L' = 3.5
    'M' = 2.8
    'N' = 2.3
    'O' = 3.6
    'P' = 2.9
    'R' = 2.025
    'S' = 1.825
    'U' = 2.025
    'V' = 1.825
    'W' = 1.3
    'X' = -1
    'Z' = 1.025
    'AA' = -1
    'AB' = -1
    'AC' = 0.825
}

# Row 99 (id 97) — corrected match data (was swapped with row 100)
Set-Row 99 @{
    'B' = 5400064
    'F' = 'Sydney FC'
    'G' = 'Newcastle Jets'
    'H' = 2
    'I' = 0
    'J' = 'H'
    'K' = 1.65
    'L' = 4
    'M' = 4.5
    'N' = 1.533
    'O' = 4.75
    'P' = 5.5
    'Q' = -1.25
    'T' = 3.5
    'U' = 2
    'V' = 1.85
    'W' = 0.5329999999999999
    'Y' = -1
    'Z' = 1.025
    'AA' = -1
    'AC' = 0.8500000000000001
}

# Row 100 (id 98) — corrected match data (was swapped with row 99)
Set-Row 100 @{
    'B' = 5404735
    'F' = 'Macarthur FC'
    'G' = 'Wellington Phoenix'
    'H' = 0
    'I' = 1
    'J' = 'A'
    'K' = 3.6
    'L' = 3.75
    'M' = 1.909
    'N' = 4
    'O' = 4
    'P' = 1.833
    'Q' = 0.5
    'T' = 3.25
    'U' = 1.925
    'V' = 1.925
    'W' = -1
    'Y' = 0.833
    'Z' = -1
    'AA' = 0.825
    'AC' = 0.925
}

# Row 214 (id 212) — updated fixture/odds data
Set-Row 214 @{
    'B' = 7126789
    'E' = 45346.14583333334
    'F' = 'Sydney FC'
    'G' = 'Melbourne City'
    'K' = 1.833
    'L' = 4
    'M' = 3.6
    'N' = 1.833
    'O' = 4.2
    'P' = 3.75
    'R' = 1.85
    'S' = 2.05
    'U' = 1.85
    'V' = 2
}

# Row 215 (id 213) — updated fixture/odds data
Set-Row 215 @{
    'B' = 7127377
    'E' = 45346.23958333334
    'F' = 'Adelaide United'
    'G' = 'Western Sydney Wanderers'
    'K' = 2.8
    'L' = 3.6
    'M' = 2.25
    'N' = 2.4
    'O' = 4
    'P' = 2.625
    'Q' = 0
    'R' = 1.85
    'S' = 2.05
    'U' = 1.85
    'V' = 2
}

# Row 216 (id 214) — updated fixture/odds data
Set-Row 216 @{
    'B' = 7127378
    'E' = 45346.32291666666
    'F' = 'Perth Glory'
    'G' = 'Wellington Phoenix'
    'K' = 2.6
    'L' = 3.5
    'M' = 2.5
    'N' = 2.375
    'O' = 3.5
    'P' = 2.875
    'Q' = -0.25
    'R' = 2.08
    'S' = 1.82
    'T' = 3
}

# Row 217 (id 215) — updated fixture/odds data
Set-Row 217 @{
    'B' = 7127376
    'E' = 45347.125
    'F' = 'Newcastle Jets'
    'G' = 'Macarthur FC'
    'K' = 2.375
    'L' = 3.6
    'M' = 2.7
    'N' = 1.95
    'O' = 4
    'P' = 3.4
    'Q' = -0.5
    'R' = 2
    'S' = 1.9
    'T' = 3.25
    'U' = 1.875
    'V' = 1.975
}

# Row 218 (id 216) — updated fixture/odds data
Set-Row 218 @{
    'B' = 7127379
    'F' = 'Melbourne Victory'
    'G' = 'Central Coast Mariners'
    'K' = 2
    'L' = 3.6
    'M' = 3.4
    'N' = 1.95
    'O' = 3.6
    'P' = 3.8
    'Q' = -0.5
    'R' = 1.92
    'S' = 1.98
    'T' = 2.75
    'U' = 1.95
    'V' = 1.9
}

# Row 219 (old id 217) is dropped entirely — the fixture now lives in row 218's data.
$ws.Rows.Item(219).Delete()
